# Updates the cryptos price/volume table (rows 2-51) to match the latest
# scrape: refreshed Price (D) and Volume(1h) (E) figures for every coin,
# plus the Polygon/OKB rows (10/11) swapping rank order.
#
# Price cells are forced to text ("@" number format) before the value is
# assigned so strings like "1.005" or "13.20" are not reinterpreted as
# numbers (which would also silently drop trailing zeros / change the
# stored value). The style is then reset to "Normal" so no extra cell
# style gets introduced versus the original workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.220.64'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.71%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.801.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.84%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.25%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.19%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.15%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5256'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3818'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.33%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07983'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.22%  '

$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.85%  '

$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.099'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.62%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.308'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.39%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.005'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.58'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.75%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.809.87'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.64%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.313'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.51%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.09'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.47%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001093'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.44%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06597'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.63%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.004'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.22%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.90%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.966'
$ws.Range('D22').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.280.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.60%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.26%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.268'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.22%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.17%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.31%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.009.76'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.87%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.355'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.95%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.09'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.28%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1082'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.64%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.056'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.42%  '

$ws.Range('E33').Value = '  +0.98%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.547'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.09%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07214'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.92%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.33%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02308'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.82%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2143'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.51%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.090'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.49%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.602'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.56%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6188'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.53%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.168'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.50%  '

$ws.Range('E43').Value = '  -2.11%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.20'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.77%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6002'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.771'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.76%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '127.18'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.86%  '

$ws.Range('E48').Value = '  +2.99%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.923'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.58%  '

$ws.Range('E50').Value = '  -1.16%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.48%  '
